$d = $word.ActiveDocument

# Locate the "Data, Technology and Strategy Consulting" paragraph that sits
# directly under the Siege Analytics PARTNER heading, and insert the three
# new bullet paragraphs right after it (before the existing bullet list).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Data, Technology and Strategy Consulting`r") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Anchor paragraph 'Data, Technology and Strategy Consulting' not found"
}

$bullet = [char]0x2022
$newText = "$bullet Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters`r" + `
           "$bullet Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States`r" + `
           "$bullet Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis`r"

$endPos = $anchor.Range.End
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertAfter($newText)
